$d = $word.ActiveDocument

$pairs = @(
    @{old = "237÷8=29, 5";   new = "194÷6=32, 2"},
    @{old = "822÷4=205, 2";  new = "138÷9=15, 3"},
    @{old = "265÷7=37, 6";   new = "492÷4=123, 0"},
    @{old = "418÷5=83, 3";   new = "292÷4=73, 0"},
    @{old = "206÷3=68, 2";   new = "406÷9=45, 1"},
    @{old = "309÷6=51, 3";   new = "457÷4=114, 1"},
    @{old = "634÷2=317, 0";  new = "449÷6=74, 5"},
    @{old = "648÷6=108, 0";  new = "354÷2=177, 0"},
    @{old = "149÷9=16, 5";   new = "868÷8=108, 4"},
    @{old = "782÷3=260, 2";  new = "747÷8=93, 3"},
    @{old = "286÷2=143, 0";  new = "533÷5=106, 3"},
    @{old = "908÷8=113, 4";  new = "266÷5=53, 1"},
    @{old = "811÷7=115, 6";  new = "880÷8=110, 0"},
    @{old = "850÷5=170, 0";  new = "315÷2=157, 1"},
    @{old = "750÷9=83, 3";   new = "434÷3=144, 2"},
    @{old = "399÷3=133, 0";  new = "164÷9=18, 2"},
    @{old = "311÷6=51, 5";   new = "612÷7=87, 3"},
    @{old = "657÷2=328, 1";  new = "272÷3=90, 2"},
    @{old = "410÷9=45, 5";   new = "167÷5=33, 2"},
    @{old = "706÷9=78, 4";   new = "161÷6=26, 5"},
    @{old = "259÷2=129, 1";  new = "366÷6=61, 0"},
    @{old = "458÷6=76, 2";   new = "362÷7=51, 5"},
    @{old = "838÷9=93, 1";   new = "840÷8=105, 0"},
    @{old = "368÷5=73, 3";   new = "566÷4=141, 2"},
    @{old = "202÷4=50, 2";   new = "602÷3=200, 2"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}

$d.Save()
